# Generate Report for Archive
# The file "3ba89fba-3641-4a12-a870-b096bd0b1b85.md" row moved from
# "Ready for handoff" to "In Translation" status in each sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 5 corresponds to 3ba89fba-...md (zh-cn in col E, de-de in col F)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"

# zh-cn sheet: row 5, column C = Status
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"

# de-de sheet: row 5, column C = Status
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
